$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.014.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.412.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -7.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "275.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3669"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3120"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.035"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06500"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.488"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.191"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001021"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.412.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05691"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -13.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -14.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.634"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.262"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.021.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.272"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "133.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.570.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.94%  "

$ws.Range("E30").Value = "  -5.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.959"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -18.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.307"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8258"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07683"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.422"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.484"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05920"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.910"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.001"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02078"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1908"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.094"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5320"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.534"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5204"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "116.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.770"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.20%  "

$ws.Range("E50").Value = "  -10.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
